# Insert a new worksheet "may18" before the first sheet (jun18), matching the
# author's workflow of prepending a new month's score card ahead of the
# existing jun18 / jul18 sheets.
$wb = $excel.ActiveWorkbook
$first = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($first)
$ws.Name = "may18"

# Column A width (matches the "bestFit" auto width Excel applies once the
# longer "Hole 10".."Hole 18" labels are present).
$ws.Columns.Item(1).ColumnWidth = 9

# Header row
$ws.Cells.Item(1,1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(1,1).Value = 43246
$ws.Cells.Item(1,2).Value = "Score"
$ws.Cells.Item(1,3).Value = "Fairway"
$ws.Cells.Item(1,4).Value = "GIR"
$ws.Cells.Item(1,5).Value = "Putts"
$ws.Cells.Item(1,6).Value = "Comment"

$ws.Cells.Item(2,1).Value = "Hole 1"
$ws.Cells.Item(2,2).Value = 4
$ws.Cells.Item(2,3).Value = "R"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(3,1).Value = "Hole 2"
$ws.Cells.Item(3,2).Value = 4
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(4,1).Value = "Hole 3"
$ws.Cells.Item(4,2).Value = 4
$ws.Cells.Item(4,3).Value = "R"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(5,1).Value = "Hole 4"
$ws.Cells.Item(5,2).Value = 7
$ws.Cells.Item(5,3).Value = "L"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(6,1).Value = "Hole 5"
$ws.Cells.Item(6,2).Value = 3
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(7,1).Value = "Hole 6"
$ws.Cells.Item(7,2).Value = 5
$ws.Cells.Item(7,3).Value = "R"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(8,1).Value = "Hole 7"
$ws.Cells.Item(8,2).Value = 4
$ws.Cells.Item(8,3).Value = "L"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(9,1).Value = "Hole 8"
$ws.Cells.Item(9,2).Value = 3
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(10,1).Value = "Hole 9"
$ws.Cells.Item(10,2).Value = 5
$ws.Cells.Item(10,3).Value = "L"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(11,1).Value = "Hole 10"
$ws.Cells.Item(11,2).Value = 5
$ws.Cells.Item(11,3).Value = "R"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(12,1).Value = "Hole 11"
$ws.Cells.Item(12,2).Value = 4
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(13,1).Value = "Hole 12"
$ws.Cells.Item(13,2).Value = 5
$ws.Cells.Item(13,3).Value = "S"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(14,1).Value = "Hole 13"
$ws.Cells.Item(14,2).Value = 4
$ws.Cells.Item(14,3).Value = "S"
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(15,1).Value = "Hole 14"
$ws.Cells.Item(15,2).Value = 5
$ws.Cells.Item(15,3).Value = "S"
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(16,1).Value = "Hole 15"
$ws.Cells.Item(16,2).Value = 3
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(17,1).Value = "Hole 16"
$ws.Cells.Item(17,2).Value = 5
$ws.Cells.Item(17,3).Value = "S"
$ws.Cells.Item(17,5).Value = 1
$ws.Cells.Item(18,1).Value = "Hole 17"
$ws.Cells.Item(18,2).Value = 6
$ws.Cells.Item(18,3).Value = "S"
$ws.Cells.Item(18,5).Value = 2
$ws.Cells.Item(19,1).Value = "Hole 18"
$ws.Cells.Item(19,2).Value = 4
$ws.Cells.Item(19,3).Value = "L"
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,6).Value = "Check this"

# Totals row
$ws.Cells.Item(20,2).Formula = "=SUM(B2:B19)"
$ws.Cells.Item(20,5).Formula = "=SUM(E2:E19)"

# Trailing blank row keeps the date-style format below the table (mirrors the
# empty, date-formatted A22 cell left over from the previous sheets).
$ws.Cells.Item(22,1).NumberFormat = "d-mmm-yy"

# Make the new sheet the active tab with the same selection the author left
# it in.
[void]$ws.Select()
[void]$ws.Range("I10").Select()
